$d = $word.ActiveDocument

# Update the date heading in the first paragraph (wdReplaceOne = 1 scopes the
# replacement to a single match; the engine applies wdReplaceAll document-wide
# regardless of the search Range, so wdReplaceOne is used everywhere below).
$d.Content.Find.Execute("2025-10-06 Monday", $true, $false, $false, $false, $false, $true, 0, $false, "2025-10-07 Tuesday", 1) | Out-Null

# Update each arithmetic problem cell in the table, addressed by (row, col)
# so that duplicate old values (e.g. "91-72=") are each replaced independently
# with their own new value.
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.Find.Execute("1+77=", $true, $false, $false, $false, $false, $true, 0, $false, "83-68=", 1) | Out-Null
$cell = $t.Cell(1, 2)
$cell.Range.Find.Execute("84-64=", $true, $false, $false, $false, $false, $true, 0, $false, "33+64=", 1) | Out-Null
$cell = $t.Cell(1, 3)
$cell.Range.Find.Execute("93+2=", $true, $false, $false, $false, $false, $true, 0, $false, "0+73=", 1) | Out-Null
$cell = $t.Cell(1, 4)
$cell.Range.Find.Execute("88-82=", $true, $false, $false, $false, $false, $true, 0, $false, "81-60=", 1) | Out-Null
$cell = $t.Cell(1, 5)
$cell.Range.Find.Execute("52-12=", $true, $false, $false, $false, $false, $true, 0, $false, "8+82=", 1) | Out-Null
$cell = $t.Cell(2, 1)
$cell.Range.Find.Execute("90-74=", $true, $false, $false, $false, $false, $true, 0, $false, "36+55=", 1) | Out-Null
$cell = $t.Cell(2, 2)
$cell.Range.Find.Execute("52+12=", $true, $false, $false, $false, $false, $true, 0, $false, "44+16=", 1) | Out-Null
$cell = $t.Cell(2, 3)
$cell.Range.Find.Execute("86+10=", $true, $false, $false, $false, $false, $true, 0, $false, "29+23=", 1) | Out-Null
$cell = $t.Cell(2, 4)
$cell.Range.Find.Execute("80-67=", $true, $false, $false, $false, $false, $true, 0, $false, "72-13=", 1) | Out-Null
$cell = $t.Cell(2, 5)
$cell.Range.Find.Execute("3+94=", $true, $false, $false, $false, $false, $true, 0, $false, "71-48=", 1) | Out-Null
$cell = $t.Cell(3, 1)
$cell.Range.Find.Execute("23-20=", $true, $false, $false, $false, $false, $true, 0, $false, "35+43=", 1) | Out-Null
$cell = $t.Cell(3, 2)
$cell.Range.Find.Execute("67-26=", $true, $false, $false, $false, $false, $true, 0, $false, "77-26=", 1) | Out-Null
$cell = $t.Cell(3, 3)
$cell.Range.Find.Execute("36+22=", $true, $false, $false, $false, $false, $true, 0, $false, "3+75=", 1) | Out-Null
$cell = $t.Cell(3, 4)
$cell.Range.Find.Execute("77-75=", $true, $false, $false, $false, $false, $true, 0, $false, "22+16=", 1) | Out-Null
$cell = $t.Cell(3, 5)
$cell.Range.Find.Execute("78+2=", $true, $false, $false, $false, $false, $true, 0, $false, "88-77=", 1) | Out-Null
$cell = $t.Cell(4, 1)
$cell.Range.Find.Execute("70+8=", $true, $false, $false, $false, $false, $true, 0, $false, "42+38=", 1) | Out-Null
$cell = $t.Cell(4, 2)
$cell.Range.Find.Execute("91-72=", $true, $false, $false, $false, $false, $true, 0, $false, "7+31=", 1) | Out-Null
$cell = $t.Cell(4, 3)
$cell.Range.Find.Execute("91-72=", $true, $false, $false, $false, $false, $true, 0, $false, "48-17=", 1) | Out-Null
$cell = $t.Cell(4, 4)
$cell.Range.Find.Execute("94-41=", $true, $false, $false, $false, $false, $true, 0, $false, "72-4=", 1) | Out-Null
$cell = $t.Cell(4, 5)
$cell.Range.Find.Execute("23+74=", $true, $false, $false, $false, $false, $true, 0, $false, "18+70=", 1) | Out-Null
$cell = $t.Cell(5, 1)
$cell.Range.Find.Execute("39+0=", $true, $false, $false, $false, $false, $true, 0, $false, "64-44=", 1) | Out-Null
$cell = $t.Cell(5, 2)
$cell.Range.Find.Execute("13+31=", $true, $false, $false, $false, $false, $true, 0, $false, "77+15=", 1) | Out-Null
$cell = $t.Cell(5, 3)
$cell.Range.Find.Execute("74-3=", $true, $false, $false, $false, $false, $true, 0, $false, "30-25=", 1) | Out-Null
$cell = $t.Cell(5, 4)
$cell.Range.Find.Execute("91-71=", $true, $false, $false, $false, $false, $true, 0, $false, "1+42=", 1) | Out-Null
$cell = $t.Cell(5, 5)
$cell.Range.Find.Execute("55+14=", $true, $false, $false, $false, $false, $true, 0, $false, "48-9=", 1) | Out-Null
$cell = $t.Cell(6, 1)
$cell.Range.Find.Execute("70+19=", $true, $false, $false, $false, $false, $true, 0, $false, "39-5=", 1) | Out-Null
$cell = $t.Cell(6, 2)
$cell.Range.Find.Execute("9+50=", $true, $false, $false, $false, $false, $true, 0, $false, "71+17=", 1) | Out-Null
$cell = $t.Cell(6, 3)
$cell.Range.Find.Execute("25-23=", $true, $false, $false, $false, $false, $true, 0, $false, "33+39=", 1) | Out-Null
$cell = $t.Cell(6, 4)
$cell.Range.Find.Execute("1+51=", $true, $false, $false, $false, $false, $true, 0, $false, "14-14=", 1) | Out-Null
$cell = $t.Cell(6, 5)
$cell.Range.Find.Execute("96-50=", $true, $false, $false, $false, $false, $true, 0, $false, "90-39=", 1) | Out-Null
$cell = $t.Cell(7, 1)
$cell.Range.Find.Execute("36+19=", $true, $false, $false, $false, $false, $true, 0, $false, "30+35=", 1) | Out-Null
$cell = $t.Cell(7, 2)
$cell.Range.Find.Execute("20+9=", $true, $false, $false, $false, $false, $true, 0, $false, "63-20=", 1) | Out-Null
$cell = $t.Cell(7, 3)
$cell.Range.Find.Execute("30+26=", $true, $false, $false, $false, $false, $true, 0, $false, "13+6=", 1) | Out-Null
$cell = $t.Cell(7, 4)
$cell.Range.Find.Execute("19-17=", $true, $false, $false, $false, $false, $true, 0, $false, "72-16=", 1) | Out-Null
$cell = $t.Cell(7, 5)
$cell.Range.Find.Execute("41+52=", $true, $false, $false, $false, $false, $true, 0, $false, "77-5=", 1) | Out-Null
$cell = $t.Cell(8, 1)
$cell.Range.Find.Execute("56-19=", $true, $false, $false, $false, $false, $true, 0, $false, "75-25=", 1) | Out-Null
$cell = $t.Cell(8, 2)
$cell.Range.Find.Execute("37+21=", $true, $false, $false, $false, $false, $true, 0, $false, "60-30=", 1) | Out-Null
$cell = $t.Cell(8, 3)
$cell.Range.Find.Execute("45+7=", $true, $false, $false, $false, $false, $true, 0, $false, "80-75=", 1) | Out-Null
$cell = $t.Cell(8, 4)
$cell.Range.Find.Execute("65-54=", $true, $false, $false, $false, $false, $true, 0, $false, "5+26=", 1) | Out-Null
$cell = $t.Cell(8, 5)
$cell.Range.Find.Execute("47+28=", $true, $false, $false, $false, $false, $true, 0, $false, "82-29=", 1) | Out-Null
$cell = $t.Cell(9, 1)
$cell.Range.Find.Execute("82-66=", $true, $false, $false, $false, $false, $true, 0, $false, "83-2=", 1) | Out-Null
$cell = $t.Cell(9, 2)
$cell.Range.Find.Execute("23+0=", $true, $false, $false, $false, $false, $true, 0, $false, "61-41=", 1) | Out-Null
$cell = $t.Cell(9, 3)
$cell.Range.Find.Execute("99-90=", $true, $false, $false, $false, $false, $true, 0, $false, "18+51=", 1) | Out-Null
$cell = $t.Cell(9, 4)
$cell.Range.Find.Execute("9+37=", $true, $false, $false, $false, $false, $true, 0, $false, "89-5=", 1) | Out-Null
$cell = $t.Cell(9, 5)
$cell.Range.Find.Execute("40-0=", $true, $false, $false, $false, $false, $true, 0, $false, "82-65=", 1) | Out-Null
$cell = $t.Cell(10, 1)
$cell.Range.Find.Execute("69-20=", $true, $false, $false, $false, $false, $true, 0, $false, "0+49=", 1) | Out-Null
$cell = $t.Cell(10, 2)
$cell.Range.Find.Execute("55+12=", $true, $false, $false, $false, $false, $true, 0, $false, "52-38=", 1) | Out-Null
$cell = $t.Cell(10, 3)
$cell.Range.Find.Execute("37+44=", $true, $false, $false, $false, $false, $true, 0, $false, "42-16=", 1) | Out-Null
$cell = $t.Cell(10, 4)
$cell.Range.Find.Execute("53+8=", $true, $false, $false, $false, $false, $true, 0, $false, "46-40=", 1) | Out-Null
$cell = $t.Cell(10, 5)
$cell.Range.Find.Execute("13+27=", $true, $false, $false, $false, $false, $true, 0, $false, "91-26=", 1) | Out-Null
$cell = $t.Cell(11, 1)
$cell.Range.Find.Execute("54-34=", $true, $false, $false, $false, $false, $true, 0, $false, "64+13=", 1) | Out-Null
$cell = $t.Cell(11, 2)
$cell.Range.Find.Execute("20+53=", $true, $false, $false, $false, $false, $true, 0, $false, "48-46=", 1) | Out-Null
$cell = $t.Cell(11, 3)
$cell.Range.Find.Execute("96-41=", $true, $false, $false, $false, $false, $true, 0, $false, "59-57=", 1) | Out-Null
$cell = $t.Cell(11, 4)
$cell.Range.Find.Execute("9+48=", $true, $false, $false, $false, $false, $true, 0, $false, "75-61=", 1) | Out-Null
$cell = $t.Cell(11, 5)
$cell.Range.Find.Execute("18-16=", $true, $false, $false, $false, $false, $true, 0, $false, "45-40=", 1) | Out-Null
$cell = $t.Cell(12, 1)
$cell.Range.Find.Execute("0+92=", $true, $false, $false, $false, $false, $true, 0, $false, "3+1=", 1) | Out-Null
$cell = $t.Cell(12, 2)
$cell.Range.Find.Execute("66-38=", $true, $false, $false, $false, $false, $true, 0, $false, "51-49=", 1) | Out-Null
$cell = $t.Cell(12, 3)
$cell.Range.Find.Execute("12+72=", $true, $false, $false, $false, $false, $true, 0, $false, "96-87=", 1) | Out-Null
$cell = $t.Cell(12, 4)
$cell.Range.Find.Execute("42-29=", $true, $false, $false, $false, $false, $true, 0, $false, "77-46=", 1) | Out-Null
$cell = $t.Cell(12, 5)
$cell.Range.Find.Execute("13+35=", $true, $false, $false, $false, $false, $true, 0, $false, "22-17=", 1) | Out-Null
$cell = $t.Cell(13, 1)
$cell.Range.Find.Execute("15+16=", $true, $false, $false, $false, $false, $true, 0, $false, "88-78=", 1) | Out-Null
$cell = $t.Cell(13, 2)
$cell.Range.Find.Execute("16+49=", $true, $false, $false, $false, $false, $true, 0, $false, "34+52=", 1) | Out-Null
$cell = $t.Cell(13, 3)
$cell.Range.Find.Execute("29+12=", $true, $false, $false, $false, $false, $true, 0, $false, "9+26=", 1) | Out-Null
$cell = $t.Cell(13, 4)
$cell.Range.Find.Execute("49-42=", $true, $false, $false, $false, $false, $true, 0, $false, "93-25=", 1) | Out-Null
$cell = $t.Cell(13, 5)
$cell.Range.Find.Execute("75-64=", $true, $false, $false, $false, $false, $true, 0, $false, "97-95=", 1) | Out-Null
$cell = $t.Cell(14, 1)
$cell.Range.Find.Execute("23+45=", $true, $false, $false, $false, $false, $true, 0, $false, "13+24=", 1) | Out-Null
$cell = $t.Cell(14, 2)
$cell.Range.Find.Execute("9+15=", $true, $false, $false, $false, $false, $true, 0, $false, "4+3=", 1) | Out-Null
$cell = $t.Cell(14, 3)
$cell.Range.Find.Execute("57-38=", $true, $false, $false, $false, $false, $true, 0, $false, "57+0=", 1) | Out-Null
$cell = $t.Cell(14, 4)
$cell.Range.Find.Execute("85+0=", $true, $false, $false, $false, $false, $true, 0, $false, "26+34=", 1) | Out-Null
$cell = $t.Cell(14, 5)
$cell.Range.Find.Execute("92-91=", $true, $false, $false, $false, $false, $true, 0, $false, "72+2=", 1) | Out-Null
$cell = $t.Cell(15, 1)
$cell.Range.Find.Execute("80-26=", $true, $false, $false, $false, $false, $true, 0, $false, "44+7=", 1) | Out-Null
$cell = $t.Cell(15, 2)
$cell.Range.Find.Execute("68-4=", $true, $false, $false, $false, $false, $true, 0, $false, "57+38=", 1) | Out-Null
$cell = $t.Cell(15, 3)
$cell.Range.Find.Execute("37+19=", $true, $false, $false, $false, $false, $true, 0, $false, "78-18=", 1) | Out-Null
$cell = $t.Cell(15, 4)
$cell.Range.Find.Execute("91+4=", $true, $false, $false, $false, $false, $true, 0, $false, "57-54=", 1) | Out-Null
$cell = $t.Cell(15, 5)
$cell.Range.Find.Execute("66+26=", $true, $false, $false, $false, $false, $true, 0, $false, "15-13=", 1) | Out-Null
$cell = $t.Cell(16, 1)
$cell.Range.Find.Execute("80-46=", $true, $false, $false, $false, $false, $true, 0, $false, "98-35=", 1) | Out-Null
$cell = $t.Cell(16, 2)
$cell.Range.Find.Execute("64+11=", $true, $false, $false, $false, $false, $true, 0, $false, "0+8=", 1) | Out-Null
$cell = $t.Cell(16, 3)
$cell.Range.Find.Execute("74-33=", $true, $false, $false, $false, $false, $true, 0, $false, "85-48=", 1) | Out-Null
$cell = $t.Cell(16, 4)
$cell.Range.Find.Execute("79-55=", $true, $false, $false, $false, $false, $true, 0, $false, "64+15=", 1) | Out-Null
$cell = $t.Cell(16, 5)
$cell.Range.Find.Execute("5+13=", $true, $false, $false, $false, $false, $true, 0, $false, "6+53=", 1) | Out-Null
$cell = $t.Cell(17, 1)
$cell.Range.Find.Execute("79-70=", $true, $false, $false, $false, $false, $true, 0, $false, "94-92=", 1) | Out-Null
$cell = $t.Cell(17, 2)
$cell.Range.Find.Execute("33+57=", $true, $false, $false, $false, $false, $true, 0, $false, "80-27=", 1) | Out-Null
$cell = $t.Cell(17, 3)
$cell.Range.Find.Execute("83-39=", $true, $false, $false, $false, $false, $true, 0, $false, "23+65=", 1) | Out-Null
$cell = $t.Cell(17, 4)
$cell.Range.Find.Execute("89-33=", $true, $false, $false, $false, $false, $true, 0, $false, "97-71=", 1) | Out-Null
$cell = $t.Cell(17, 5)
$cell.Range.Find.Execute("0+44=", $true, $false, $false, $false, $false, $true, 0, $false, "36+53=", 1) | Out-Null
$cell = $t.Cell(18, 1)
$cell.Range.Find.Execute("70-49=", $true, $false, $false, $false, $false, $true, 0, $false, "53-17=", 1) | Out-Null
$cell = $t.Cell(18, 2)
$cell.Range.Find.Execute("73+10=", $true, $false, $false, $false, $false, $true, 0, $false, "12-1=", 1) | Out-Null
$cell = $t.Cell(18, 3)
$cell.Range.Find.Execute("73-13=", $true, $false, $false, $false, $false, $true, 0, $false, "10+40=", 1) | Out-Null
$cell = $t.Cell(18, 4)
$cell.Range.Find.Execute("18+80=", $true, $false, $false, $false, $false, $true, 0, $false, "32-18=", 1) | Out-Null
$cell = $t.Cell(18, 5)
$cell.Range.Find.Execute("18-11=", $true, $false, $false, $false, $false, $true, 0, $false, "5+76=", 1) | Out-Null
$cell = $t.Cell(19, 1)
$cell.Range.Find.Execute("7+12=", $true, $false, $false, $false, $false, $true, 0, $false, "89-11=", 1) | Out-Null
$cell = $t.Cell(19, 2)
$cell.Range.Find.Execute("40+31=", $true, $false, $false, $false, $false, $true, 0, $false, "68-5=", 1) | Out-Null
$cell = $t.Cell(19, 3)
$cell.Range.Find.Execute("79-14=", $true, $false, $false, $false, $false, $true, 0, $false, "21-1=", 1) | Out-Null
$cell = $t.Cell(19, 4)
$cell.Range.Find.Execute("42-40=", $true, $false, $false, $false, $false, $true, 0, $false, "2+50=", 1) | Out-Null
$cell = $t.Cell(19, 5)
$cell.Range.Find.Execute("54+24=", $true, $false, $false, $false, $false, $true, 0, $false, "54+43=", 1) | Out-Null
$cell = $t.Cell(20, 1)
$cell.Range.Find.Execute("93-83=", $true, $false, $false, $false, $false, $true, 0, $false, "24+50=", 1) | Out-Null
$cell = $t.Cell(20, 2)
$cell.Range.Find.Execute("92-81=", $true, $false, $false, $false, $false, $true, 0, $false, "27+39=", 1) | Out-Null
$cell = $t.Cell(20, 3)
$cell.Range.Find.Execute("14+42=", $true, $false, $false, $false, $false, $true, 0, $false, "29-29=", 1) | Out-Null
$cell = $t.Cell(20, 4)
$cell.Range.Find.Execute("93-76=", $true, $false, $false, $false, $false, $true, 0, $false, "44-43=", 1) | Out-Null
$cell = $t.Cell(20, 5)
$cell.Range.Find.Execute("50-43=", $true, $false, $false, $false, $false, $true, 0, $false, "25+10=", 1) | Out-Null
